$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.993.65"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.205.59"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.76"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.508"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.34%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "29.97"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0777"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.91"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +6.36%  "
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.44"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.551.10"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.72"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.203.09"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.724"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.907.66"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0882"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.76"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.14"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.90"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.80"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.41"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.13"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.45"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.26"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0708"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.82"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0973"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.34"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.61%  "
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.119.64"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0267"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.42"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.63"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.423.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.48"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("E51").Value = "  +0.78%  "
